# Update cryptos list values per the Nov 22 2024 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'98.947.05"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.20%  "

$ws.Range("D3").Value = "'3.385.08"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +8.25%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'262.31"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +8.75%  "

$ws.Range("D6").Value = "'633.64"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.54%  "

$ws.Range("D7").Value = "'1.39"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +24.37%  "

$ws.Range("E8").Value = "  +2.19%  "

$ws.Range("E9").Value = "  -0.05%  "

$ws.Range("D10").Value = "'0.884"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +12.58%  "

$ws.Range("D11").Value = "'3.383.91"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +8.33%  "

$ws.Range("E12").Value = "  +1.41%  "

$ws.Range("D13").Value = "'98.520.87"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.44%  "

$ws.Range("D14").Value = "'36.28"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +6.81%  "

$ws.Range("D15").Value = "'0.0000249"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.87%  "

$ws.Range("D16").Value = "'4.006.01"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +7.98%  "

$ws.Range("E17").Value = "  +3.29%  "

$ws.Range("D18").Value = "'3.386.82"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +8.35%  "

$ws.Range("D19").Value = "'3.61"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.32%  "

$ws.Range("D20").Value = "'15.27"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.35%  "

$ws.Range("D21").Value = "'496.03"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.53%  "

$ws.Range("D22").Value = "'6.21"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +9.03%  "

$ws.Range("E23").Value = "  +9.27%  "

$ws.Range("D24").Value = "'9.41"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +6.17%  "

$ws.Range("D25").Value = "'5.79"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.54%  "

$ws.Range("D26").Value = "'90.43"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.55%  "

$ws.Range("D27").Value = "'12.10"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.58%  "

$ws.Range("D28").Value = "'3.559.63"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +8.19%  "

$ws.Range("D29").Value = "'0.282"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +18.46%  "

$ws.Range("D30").Value = "'0.200"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +14.11%  "

$ws.Range("D31").Value = "'0.995"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.35%  "

$ws.Range("D32").Value = "'0.133"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +6.10%  "

$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").Value = "'1.00"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +19.07%  "

$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'9.60"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +5.85%  "

$ws.Range("D35").Value = "'27.94"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +5.13%  "

$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.151"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.57%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D37").Value = "'7.36"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.07%  "

$ws.Range("D38").Value = "'1.98"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +5.74%  "

$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").Value = "'0.473"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +7.76%  "

$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "'506.92"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.85%  "

$ws.Range("D41").Value = "'24.85"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.40%  "

$ws.Range("E42").Value = "  +2.51%  "

$ws.Range("D43").Value = "'3.71"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.33%  "

$ws.Range("D44").Value = "'3.37"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +5.16%  "

$ws.Range("E45").Value = "  +13.07%  "

$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("D47").Value = "'161.68"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.87%  "

$ws.Range("D48").Value = "'1.96"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.33%  "

$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.845"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +15.11%  "

$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").Value = "'4.69"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +6.67%  "

$ws.Range("D51").Value = "'46.45"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +4.37%  "
